$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('I1').Value = 'Other found locations'

$ws.Range('I2').Value = ''

$ws.Range('F3').Value = 'not found'
$ws.Range('G3').Value = 'N/A'
$ws.Range('E3').Value = '[Chai%X%coreGivesNoEmail%2,   Deliang%Huang%coreGivesNoEmail%2,   Hong%Yu%coreGivesNoEmail%2,   Jun%Chen%coreGivesNoEmail%0,   Lei%Liu%coreGivesNoEmail%0,   National%Health Commission of the People’s Republic of China%coreGivesNoEmail%2,   Pengcheng%Ou%coreGivesNoEmail%2,   Qing%He%coreGivesNoEmail%2,   Qingxian%Cai%coreGivesNoEmail%2,   Rossaint%J%coreGivesNoEmail%2,   Yang%Fu%coreGivesNoEmail%2,   Yiming%Zhang%coreGivesNoEmail%2,   Yinan%Su%coreGivesNoEmail%2,   Zhang%Xia%coreGivesNoEmail%2,   Zhenghua%Ma%coreGivesNoEmail%2,   Zhibin%Zhu%coreGivesNoEmail%2,   Zhiwei%Li%coreGivesNoEmail%2]'
$ws.Range('I3').Value = ''

$ws.Range('E4').Value = '[Yong%Gao%NULL%0,    Tuantuan%Li%NULL%0,    Mingfeng%Han%NULL%0,    Xiuyong%Li%NULL%0,    Dong%Wu%NULL%0,    Yuanhong%Xu%NULL%0,    Yulin%Zhu%NULL%0,    Yan%Liu%NULL%0,    Xiaowu%Wang%wangxiaowu19880218@126.com%0,    Linding%Wang%wanglinding@ahmu.edu.cn%0]'
$ws.Range('I4').Value = '_PMC'

$ws.Range('E5').Value = '[Zhongliang%Wang%NULL%0,    Bohan%Yang%NULL%0,    Qianwen%Li%NULL%0,    Lu%Wen%NULL%0,    Ruiguang%Zhang%zrg27@163.com%0]'
$ws.Range('I5').Value = '_PMC'

$ws.Range('E6').Value = '[Chuan%Qin%NULL%0,    Luoqi%Zhou%NULL%0,    Ziwei%Hu%NULL%0,    Shuoqi%Zhang%NULL%0,    Sheng%Yang%NULL%0,    Yu%Tao%NULL%0,    Cuihong%Xie%NULL%0,    Ke%Ma%NULL%0,    Ke%Shang%NULL%0,    Wei%Wang%NULL%0,    Dai-Shi%Tian%tiands@tjh.tjmu.edu.cn%0]'
$ws.Range('I6').Value = '_PMC'

$ws.Range('I7').Value = ''

$ws.Range('I8').Value = ''

$ws.Range('E9').Value = '[Ruirui%Wang%NULL%0,    Min%Pan%NULL%0,    Xiumei%Zhang%NULL%0,    Mingfeng%Han%fyhmf@163.com%0,    Xiaoyun%Fan%13956988552@126.com%0,    Fengde%Zhao%NULL%0,    Manli%Miao%NULL%0,    Jing%Xu%NULL%0,    Minglong%Guan%NULL%0,    Xia%Deng%NULL%0,    Xu%Chen%NULL%0,    Leilei%Shen%NULL%0]'
$ws.Range('I9').Value = '_PMC_elsevier'

$ws.Range('E10').Value = '[Zhe%Zhu%NULL%1,    Ting%Cai%NULL%1,    Lingyan%Fan%NULL%1,    Kehong%Lou%NULL%1,    Xin%Hua%NULL%1,    Zuoan%Huang%NULL%1,    Guosheng%Gao%NULL%1]'
$ws.Range('I10').Value = '_PMC_elsevier'

$ws.Range('E11').Value = '[Xiaohua%Chen%NULL%1,    Binghong%Zhao%NULL%1,    Yueming%Qu%NULL%1,    Yurou%Chen%NULL%1,    Jie%Xiong%NULL%0,    Yong%Feng%NULL%0,    Dong%Men%NULL%1,    Qianchuan%Huang%NULL%1,    Ying%Liu%NULL%0,    Bo%Yang%NULL%0,    Jinya%Ding%NULL%1,    Feng%Li%fli222@whu.edu.cn%0]'
$ws.Range('I11').Value = '_PMC'

$ws.Range('E12').Value = '[Ming%Ding%NULL%1,    Qiang%Zhang%NULL%1,    Qing%Li%NULL%1,    Ting%Wu%NULL%1,    Ying-zi%Huang%NULL%1]'
$ws.Range('I12').Value = '_PMC_elsevier'

$ws.Range('E13').Value = '[Peder L.%Myhre%p.l.myhre@medisin.uio.no%1,    Christian%Prebensen%christian.prebensen@gmail.com%1,    Heidi%Strand%Heidi.Strand@ahus.no%1,    Ragnhild%Røysland%ragnhiro@medisin.uio.no%1,    Christine M.%Jonassen%chrjon@so-hf.no%1,    Anbjørg%Rangberg%Anbjorg.Rangberg@so-hf.no%1,    Vibecke%Sørensen%Vibecke.Sorensen@ahus.no%1,    Signe%Søvik%signe.sovik@medisin.uio.no%1,    Helge%Røsjø%helge.rosjo@medisin.uio.no%1,    My%Svensson%m.h.s.svensson@medisin.uio.no%1,    Jan%Erik Berdal%jan-erik.berdal@hotmail.com%1,    Torbjørn%Omland%NULL%1]'
$ws.Range('I13').Value = '_PMC'

$ws.Range('E14').Value = '[Ana Carolina%Gadotti%NULL%1,    Marina%de Castro Deus%NULL%1,    Joao Paulo%Telles%NULL%1,    Rafael%Wind%NULL%1,    Marina%Goes%NULL%1,    Roberta%Garcia Charello Ossoski%NULL%1,    Alessandra Michalski%de Padua%NULL%1,    Lucia%de Noronha%NULL%1,    Andrea%Moreno-Amaral%NULL%1,    Cristina Pellegrino%Baena%NULL%1,    Felipe Francisco%Tuon%NULL%1]'
$ws.Range('I14').Value = '_PMC_elsevier'

$ws.Range('E15').Value = '[Rocio%Laguna-Goya%NULL%1,    Alberto%Utrero-Rico%NULL%1,    Paloma%Talayero%NULL%1,    Maria%Lasa-Lazaro%NULL%1,    Angel%Ramirez-Fernandez%NULL%1,    Laura%Naranjo%NULL%1,    Alejandro%Segura-Tudela%NULL%1,    Oscar%Cabrera-Marante%NULL%1,    Edgar%Rodriguez de Frias%NULL%1,    Rocio%Garcia-Garcia%NULL%1,    Mario%Fernández-Ruiz%NULL%0,    Jose Maria%Aguado%NULL%1,    Joaquin%Martinez-Lopez%NULL%1,    Elena Ana%Lopez%NULL%1,    Mercedes%Catalan%NULL%1,    Antonio%Serrano%NULL%1,    Estela%Paz-Artal%NULL%1]'
$ws.Range('I15').Value = '_PMC_elsevier'

$ws.Range('E16').Value = '[Tao%Chen%NULL%0,    Di%Wu%NULL%0,    Huilong%Chen%NULL%0,    Weiming%Yan%NULL%0,    Danlei%Yang%NULL%0,    Guang%Chen%NULL%0,    Ke%Ma%NULL%0,    Dong%Xu%NULL%0,    Haijing%Yu%NULL%0,    Hongwu%Wang%NULL%0,    Tao%Wang%NULL%0,    Wei%Guo%NULL%0,    Jia%Chen%NULL%0,    Chen%Ding%NULL%0,    Xiaoping%Zhang%NULL%0,    Jiaquan%Huang%NULL%0,    Meifang%Han%NULL%0,    Shusheng%Li%NULL%0,    Xiaoping%Luo%NULL%0,    Jianping%Zhao%NULL%0,    Qin%Ning%NULL%0]'
$ws.Range('I16').Value = '_PMC'

$ws.Range('E17').Value = '[Lang%Wang%NULL%0,    Wenbo%He%NULL%0,    Xiaomei%Yu%NULL%0,    Dalong%Hu%NULL%0,    Mingwei%Bao%NULL%0,    Huafen%Liu%NULL%0,    Jiali%Zhou%NULL%0,    Hong%Jiang%NULL%0]'
$ws.Range('I17').Value = '_PMC_elsevier'

$ws.Range('E18').Value = '[Fei%Zhou%NULL%0,    Ting%Yu%NULL%0,    Ronghui%Du%NULL%0,    Guohui%Fan%NULL%0,    Ying%Liu%NULL%0,    Zhibo%Liu%NULL%0,    Jie%Xiang%NULL%0,    Yeming%Wang%NULL%0,    Bin%Song%NULL%0,    Xiaoying%Gu%NULL%0,    Lulu%Guan%NULL%0,    Yuan%Wei%NULL%0,    Hui%Li%NULL%0,    Xudong%Wu%NULL%0,    Jiuyang%Xu%NULL%0,    Shengjin%Tu%NULL%0,    Yi%Zhang%NULL%0,    Hua%Chen%NULL%0,    Bin%Cao%NULL%0]'
$ws.Range('I18').Value = '_PMC_elsevier'

$ws.Range('E19').Value = '[Haiying%Sun%NULL%0,    Ruoqi%Ning%NULL%0,    Yu%Tao%NULL%0,    Chong%Yu%NULL%0,    Xiaoyan%Deng%NULL%0,    Caili%Zhao%NULL%0,    Silu%Meng%NULL%0,    Fangxu%Tang%89650793@qq.com%0,    Dong%Xu%89650793@qq.com%0,    Dong%Xu%89650793@qq.com%0]'
$ws.Range('I19').Value = '_PMC'

$ws.Range('E20').Value = '[Junli%Fan%NULL%1,    Hui%Wang%NULL%0,    Guangming%Ye%NULL%2,    Xiaoling%Cao%NULL%1,    Xianqun%Xu%NULL%1,    Wenbin%Tan%wenbin.tan@uscmed.sc.edu%1,    Yongxi%Zhang%znact1936@126.com%0]'
$ws.Range('I20').Value = '_PMC_elsevier'

$ws.Range('E21').Value = '[Mathilda%Mandel%Mati.Mandel@sheba.health.gov.il%1,    Gil%Harari%NULL%1,    Michael%Gurevich%NULL%1,    Anat%Achiron%NULL%1]'
$ws.Range('I21').Value = '_PMC_elsevier'

$ws.Range('E22').Value = '[Jose J.%Guirao%NULL%1,    Carmen M.%Cabrera%NULL%1,    Natalia%Jiménez%NULL%1,    Laura%Rincón%NULL%1,    José M.%Urra%NULL%1]'
$ws.Range('I22').Value = '_PMC_elsevier'

$ws.Range('E23').Value = '[Chunjin%Ke%NULL%1,    Chong%Yu%NULL%0,    Daoyuan%Yue%NULL%1,    Xing%Zeng%NULL%1,    Zhiquan%Hu%NULL%1,    Chunguang%Yang%NULL%1]'
$ws.Range('I23').Value = '_PMC_elsevier'

$ws.Range('E24').Value = '[Fuyang%Chen%NULL%0,    Wenwu%Sun%NULL%0,    Shengrong%Sun%NULL%0,    Zhiyu%Li%lizhiyu@whu.edu.cn%0,    Zhong%Wang%zhongwangchn@whu.edu.cn%0,    Li%Yu%yuliwhzxyy@163.com%0,    Li%Yu%yuliwhzxyy@163.com%0]'
$ws.Range('I24').Value = '_PMC'

$ws.Range('E25').Value = '[Hui%Chen%NULL%1,    Jun%Wang%NULL%0,    Nan%Su%NULL%1,    Xiebing%Bao%NULL%1,    Yongsheng%Li%dr_ysli@126.com%1,    Jun%Jin%jinjun0514@163.com%1]'
$ws.Range('I25').Value = '_PMC_Springer'

$ws.Range('E26').Value = '[Chenze%Li%NULL%0,    Jiangang%Jiang%NULL%0,    Feng%Wang%NULL%0,    Ning%Zhou%NULL%0,    Giacomo%Veronese%NULL%0,    Javid J.%Moslehi%NULL%0,    Enrico%Ammirati%NULL%0,    Dao Wen%Wang%NULL%0]'
$ws.Range('I26').Value = '_PMC_elsevier'

$ws.Range('I27').Value = ''

$ws.Range('E28').Value = '[Zhao-Hua%Wang%NULL%1,    Chang%Shu%NULL%0,    Xiao%Ran%NULL%1,    Cui-Hong%Xie%NULL%1,    Lei%Zhang%NULL%0]'
$ws.Range('I28').Value = '_PMC'

$ws.Range('E29').Value = '[Benping%Zhang%NULL%1,    Chen%Dong%NULL%1,    Shengzhong%Li%NULL%1,    Xiaoqing%Song%NULL%1,    Wang%Wei%NULL%1,    Li%Liu%NULL%0]'
$ws.Range('I29').Value = '_PMC'

$ws.Range('E30').Value = '[Lin%Zhang%NULL%1,    Bin%Huang%NULL%1,    Hongzhen%Xia%NULL%1,    Hua%Fan%NULL%1,    Muxin%Zhu%NULL%1,    Liping%Zhu%NULL%1,    Huan%Zhang%NULL%1,    Xiaogen%Tao%NULL%1,    Shaohui%Cheng%NULL%1,    Jian%Chen%NULL%0]'
$ws.Range('I30').Value = '_PMC'
